# The deck's slide theme ("Integral") and its Notes Master theme ("Office
# Theme") get swapped: the slides end up on the plain default "Office Theme"
# palette while the Notes Master keeps the old "Integral" palette.
#
# The PowerPoint object model only exposes the 12-color DrawingML theme
# scheme for the *slide* theme (via Slide.ThemeColorScheme / the shared
# theme used by every slide through the one SlideMaster) -- there is no
# writable property for the Notes Master's own theme colors, so this
# reproduces the reachable half of the change: re-pointing every theme
# color used by the slides from the "Integral" palette to the stock
# "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
